# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: "Created functions to get season record" — the
# per-player table gains three new trailing columns (AD:AF) carrying the
# team's season Wins/Losses/Ties, repeated for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 --------------------------
# Copy formatting from the last existing header cell (AC1, style index 1:
# bold font + thin border + center/top alignment) onto the new header
# cells so they visually match the rest of the header row, then set text.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# --- Data rows (2-47): season record repeated for every player ---------
$wins = 86
$losses = 76
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-47"
